$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 updates
$ws.Range("I7").Value = 6.25
$ws.Range("K7").Value = 2.05
$ws.Range("L7").Value = 6.5
$ws.Range("Q7").Value = 2.35
$ws.Range("R7").Value = 1.57
$ws.Range("S7").Value = 3.6
$ws.Range("T7").Value = 1.28
$ws.Range("U7").Value = 4.33
$ws.Range("V7").Value = 1.2
$ws.Range("W7").Value = 1.5
$ws.Range("X7").Value = 2.5
$ws.Range("Y7").Value = 2.25
$ws.Range("Z7").Value = 1.57
$ws.Range("AR7").Value = 1.78
$ws.Range("AS7").Value = 2.1

# Row 14 updates
$ws.Range("G14").Value = 1.4
$ws.Range("H14").Value = 4.75
$ws.Range("I14").Value = 7.5
$ws.Range("J14").Value = 1.91
$ws.Range("L14").Value = 7
$ws.Range("M14").Value = 1.03
$ws.Range("N14").Value = 15
$ws.Range("Y14").Value = 1.95
$ws.Range("Z14").Value = 1.8
$ws.Range("AF14").Value = 26
$ws.Range("AH14").Value = 9
$ws.Range("AK14").Value = 301
$ws.Range("AM14").Value = 41
